# Apply the arithmetic-answer updates described in the commit diff.
# Each entry replaces the old "a+b=c" / "a-b=c" text with the new one,
# using Word's Find/Replace (exact match, whole string) against the
# document content.

$d = $word.ActiveDocument

$replacements = @(
    @("7+31=38", "8+65=73"),
    @("11+75=86", "38+46=84"),
    @("43+23=66", "95-48=47"),
    @("16+54=70", "42+23=65"),
    @("25+41=66", "65+22=87"),
    @("44-18=26", "15+21=36"),
    @("51+24=75", "10+5=15"),
    @("54+1=55", "31-23=8"),
    @("39-2=37", "54-1=53"),
    @("72-14=58", "25+35=60"),
    @("63-58=5", "91-39=52"),
    @("55-46=9", "8-1=7"),
    @("15+81=96", "93-13=80"),
    @("42+31=73", "93-62=31"),
    @("79-18=61", "6+74=80"),
    @("63+12=75", "40-7=33"),
    @("13+38=51", "81-12=69"),
    @("23+31=54", "33+30=63"),
    @("38+38=76", "75-17=58"),
    @("31+0=31", "25+15=40"),
    @("30-9=21", "1+9=10"),
    @("53-7=46", "50+31=81"),
    @("71-17=54", "21+4=25"),
    @("15-8=7", "60-18=42"),
    @("52-9=43", "70+20=90"),
    @("77-72=5", "73-44=29"),
    @("1+39=40", "7+34=41"),
    @("43+6=49", "22+53=75"),
    @("35-4=31", "47+37=84"),
    @("60-8=52", "58-38=20"),
    @("39-36=3", "8+76=84"),
    @("36-21=15", "61-32=29"),
    @("18+18=36", "38-11=27"),
    @("37+51=88", "62-40=22"),
    @("59-17=42", "28-8=20"),
    @("69+27=96", "17+34=51"),
    @("57+15=72", "84+12=96"),
    @("42-40=2", "88-37=51"),
    @("10+40=50", "37-4=33"),
    @("29-29=0", "5+91=96"),
    @("54-33=21", "38-6=32"),
    @("89-3=86", "98-87=11"),
    @("86-43=43", "25+25=50"),
    @("62-10=52", "64-53=11"),
    @("89-0=89", "21+38=59"),
    @("25+32=57", "4+57=61"),
    @("17+37=54", "2+61=63"),
    @("54+6=60", "51+41=92"),
    @("10+76=86", "57+6=63"),
    @("73+26=99", "92+0=92"),
    @("67-63=4", "66-34=32"),
    @("46+25=71", "19+64=83"),
    @("85-30=55", "30-13=17"),
    @("0+29=29", "28+45=73"),
    @("54+32=86", "98-69=29"),
    @("35-16=19", "28+44=72"),
    @("19+30=49", "42+33=75"),
    @("33+9=42", "12+72=84"),
    @("58-26=32", "75-68=7"),
    @("50+23=73", "94-84=10"),
    @("7+7=14", "7+1=8"),
    @("23+10=33", "94-16=78"),
    @("96-30=66", "91-90=1"),
    @("66+10=76", "22+0=22"),
    @("51+10=61", "40-28=12"),
    @("50-31=19", "42+30=72"),
    @("6+13=19", "75-31=44"),
    @("77-38=39", "8+50=58"),
    @("20-10=10", "63-36=27"),
    @("49-22=27", "78-51=27"),
    @("0+0=0", "27-7=20"),
    @("57-14=43", "62+22=84"),
    @("3+71=74", "62-30=32"),
    @("16+29=45", "10+49=59"),
    @("78-31=47", "97-21=76"),
    @("63-9=54", "88+7=95"),
    @("69-14=55", "37-19=18"),
    @("35-31=4", "98-45=53"),
    @("80-34=46", "89-45=44"),
    @("67-47=20", "30+24=54"),
    @("96-59=37", "18+57=75"),
    @("87-46=41", "69+19=88"),
    @("92-70=22", "98-61=37"),
    @("9+29=38", "11+23=34"),
    @("42+42=84", "10-0=10"),
    @("24+63=87", "30-29=1"),
    @("10+31=41", "42+38=80"),
    @("34-9=25", "22+14=36"),
    @("80-39=41", "86-12=74"),
    @("44+45=89", "42+36=78"),
    @("47-21=26", "73-70=3"),
    @("4+4=8", "22-4=18"),
    @("20+78=98", "7+61=68"),
    @("54-34=20", "44-29=15"),
    @("96-69=27", "57-29=28"),
    @("62-14=48", "53-43=10"),
    @("43+2=45", "33+40=73"),
    @("74-73=1", "21+55=76"),
    @("51+31=82", "49+22=71"),
    @("46-40=6", "72-41=31")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}
